$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the top. This shifts every existing row
# (including the old header row 1 and all 64 data rows) down by one,
# carrying their values AND formatting (styles) with them.
$ws.Rows.Item(1).Insert()

# At this point:
#  - row 1 is a fresh blank row (no formatting)
#  - row 2 holds what used to be row 1 (the text header labels), still
#    bearing the bold/centered/bordered header style
#  - rows 3..66 hold the old data rows 2..65

# Copy the header formatting (bold font, border, centered/top alignment)
# from row 2 onto the new row 1 before we strip it from row 2.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 1 with the 0-based column-index values.
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# Strip the header formatting back off row 2 so it looks like an
# ordinary data row again (no bold, no border, default alignment/style).
$ws.Range("A2:L2").ClearFormats()

# The K2/L2 cells ("thread_size" / "material_surface" column labels)
# are cleared out entirely in the new layout.
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
